# New names for the WSA data feeds — rename the eight "WSA*" item labels
# in column A (rows 12-19) to their new descriptive names, mark the
# renamed cells (and a helper column F used while testing) as Text
# format, and leave the selection on the first renamed cell (A12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the items. Order matters for shared-string table layout
# (row 16 is renamed before row 15 here, matching how it was done live).
$ws.Range("A12").Value = "SteelOxygenBlownConverters"
$ws.Range("A13").Value = "SteelElectricFurnaces"
$ws.Range("A14").Value = "FlatRolledProducts"
$ws.Range("A16").Value = "LongRolledProducts"
$ws.Range("A15").Value = "Ingots"
$ws.Range("A17").Value = "SteelOpenHearthFurnaces"
$ws.Range("A18").Value = "PigIron"
$ws.Range("A19").Value = "SpongeIron"

# Format a subset of the renamed cells as Text.
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A16").NumberFormat = "@"

# Also format the (otherwise empty) helper column F alongside, rows 12-19.
$ws.Range("F12:F19").NumberFormat = "@"

# Leave the selection where the edits were made.
[void]$ws.Range("A12").Select()
